# "SUPER HUGE COMMIT FROM THE SERVER"
# Replaces the single combined "Все классы" roster with the "9C класс"
# roster: renames the sheet/tab, drops the now-unused "Класс" column,
# replaces the 4 placeholder students with the real 13-student 9C list,
# and moves/updates the trailing summary row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the sheet -------------------------------------------------
$ws.Name = "9C класс"

# --- Grab the summary row's formatting (bold, right aligned) before we
#     overwrite row 8 with roster data, and stamp it onto the new
#     summary row 17. -----------------------------------------------
$ws.Range("F8").Copy()
$ws.Range("F17").PasteSpecial(-4122)

# --- Drop the "Класс" column (H) --------------------------------------
$ws.Columns.Item(8).Delete()

# --- Student roster (A:#, B:ФИО Ученика, C:ИИН, D:ФИО Родителя, ------
#     E:Номер, F:Адрес, G:Оплата) ---------------------------------------
$students = @(
  @("Айдашева Айзере Темирбековна", "110802603492", "Айдашева Гульзат Аслановна", "+7 (701) 556-54-58", "Алтын орда 19 кв 57", 35000),
  @("Базарбай Нұртас Құрманбекұлы", "100509552658", "Ешмагамбетова Айнур ", "+7 (771) 828-33-44", "Мустафа Шокай 48В к1 кв 18", 147600),
  @("Базарқұл Гүлназ Аманқұлқызы", "111043602043", "Махамбетова Жаннұр Казиевна", "+7 (708) 176-94-96", "Құрмашев 28", 25000),
  @("Батырбаева Айлин Еркиновна", "120719603413", "Батырбаева Гулсамал Есеновна", "+7 (707) 883-37-33", "Ветеран 2 д31", 81000),
  @("Жумакаева Молдир Бахадыровна", "101130602600", "Ермагамбетова Диляра ", "+7 (705) 473-60-75", "уч №143 Шабыт 124", 250000),
  @("Изтилеу Аяулым Асхатқызы", "110813603790", "сейлханова Рысгул Мирхановна", "+7 (701) 615-70-79", "м.шоқай71", 50000),
  @("Муслимова Ляйсан Дамировна", "110518603617", "Муслимова Динара Тлеповна", "+7 (777) 817-32-96", "Батыс 2 дом 7 корпус 4", 0),
  @("Назарбаева Лана Данияровна", "100914653238", "Назарбаева Асель Маратовна", "+7 (701) 920-30-48", "Саздинское лесничество 197а, кв3", 50000),
  @("Орынбасар Аянат Дулаткызы", "110125602593", "Құдайбергенова Әсел Әлибекқызы", "+7 (702) 106-56-85", "Мангилик ел 5, к4, кв10", 102000),
  @("Өтеп Әмина Нұрасханқызы", "110218604684", "Танирбергенова Айгүл Ельбаевна", "+7 (705) 474-77-84", "Жаңақоныс, Инабат 22", 43700),
  @("Пирвердиева Тамила Заургызы", "110125603660", "Пирвердиева Динара Муханбетовна", "+7 (702) 201-05-50", "Акжар - 2, Райымбек батыра 8", 50000),
  @("Шахаров Әлинур Асхатұлы", "110221500831", "Шахарова Алия Булаткалиевна", "+7 (701) 554-04-99", "Сәңкібай батыр көшесі, 28Вк3, 36 пәтер", 38000),
  @("Шаяхмет Дәулет Мирланұлы", "101112502745", "Едилова Динара Ерсаиновна", "+7 (701) 222-17-86", "Ғ.жубанова 83", 100000)
)

$firstDataRow = 2
$lastDataRow = $firstDataRow + $students.Count - 1

# Rows 6..14 don't exist yet in the sheet - stamp them with the same
# A:G formatting as row 5 before writing values into them.
$formatRow = 5
for ($r = $formatRow + 1; $r -le $lastDataRow; $r++) {
  $ws.Range("A$formatRow`:G$formatRow").Copy()
  $ws.Range("A$r`:G$r").PasteSpecial(-4122)
}

$r = $firstDataRow
$n = 1
foreach ($s in $students) {
  $ws.Cells.Item($r, 1).Value = $n
  $ws.Cells.Item($r, 2).Value = $s[0]
  $ws.Cells.Item($r, 3).Value = $s[1]
  $ws.Cells.Item($r, 4).Value = $s[2]
  $ws.Cells.Item($r, 5).Value = $s[3]
  $ws.Cells.Item($r, 6).Value = $s[4]
  $ws.Cells.Item($r, 7).Value = $s[5]
  $r = $r + 1
  $n = $n + 1
}

# The ИИН column (C) is a digit-only string (national ID number), but a
# plain numeric-looking assignment auto-coerces to a real number. Force
# it back to text, then restore the plain (border-only) cell style that
# coercion-avoidance would otherwise bump onto a freshly-minted xf.
$ws.Range("C$firstDataRow`:C$lastDataRow").NumberFormat = "@"
$r = $firstDataRow
foreach ($s in $students) {
  $ws.Cells.Item($r, 3).Value = $s[1]
  $r = $r + 1
}
$ws.Range("D$firstDataRow`:D$lastDataRow").Copy()
$ws.Range("C$firstDataRow`:C$lastDataRow").PasteSpecial(-4122)

# --- Trailing summary row now lives at row 17 --------------------------
$ws.Cells.Item(17, 6).Value = "9C класс - 13 учеников"

# --- Page margins back to Excel defaults -------------------------------
$ws.PageSetup.LeftMargin = 0.75 * 72
$ws.PageSetup.RightMargin = 0.75 * 72
$ws.PageSetup.TopMargin = 1 * 72
$ws.PageSetup.BottomMargin = 1 * 72
$ws.PageSetup.HeaderMargin = 0.5 * 72
$ws.PageSetup.FooterMargin = 0.5 * 72

Write-Host "edit complete"
